$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A104").Value = "Only 1099/W2 | Golang Devleoper@ Plano, TX"
$ws.Range("B104").Value = "https://www.dice.com/job-detail/72b40b8d-3199-4b85-9be5-730ee078ac69"
$ws.Range("C104").Value = "Plano, Texas"
$ws.Range("D104").Value = "Contract"
$ws.Range("E104").Value = "Depends on Experience"
$ws.Range("F104").Value = "InfiCare Technologies"

$ws.Range("A105").Value = "Senior GoLang Developer"
$ws.Range("B105").Value = "https://www.dice.com/job-detail/5e65c892-43e2-4daf-8978-8491e4c9af4b"
$ws.Range("C105").Value = "Hybrid in Plano, Texas"
$ws.Range("D105").Value = "Contract"
$ws.Range("E105").Value = "Depends on Experience"
$ws.Range("F105").Value = "Avtech Solutions"
